# Apply updated "dSF" (column F) values to Sheet1, per the data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F4").Value = 7
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -3
$ws.Range("F19").Value = -1
$ws.Range("F21").Value = 11
$ws.Range("F24").Value = 1
$ws.Range("F27").Value = -6
$ws.Range("F29").Value = -1
$ws.Range("F32").Value = 1
$ws.Range("F36").Value = -7
$ws.Range("F37").Value = 1
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 3
$ws.Range("F45").Value = -4
$ws.Range("F46").Value = -2
$ws.Range("F51").Value = -1
$ws.Range("F52").Value = 4
$ws.Range("F63").Value = -3
$ws.Range("F64").Value = 0
$ws.Range("F65").Value = -4
$ws.Range("F69").Value = 8
$ws.Range("F73").Value = -4

$wb.Save()
